$d = $word.ActiveDocument

$d.Content.Find.Execute("guzman", $true, $false, $false, $false, $false, $true, 1, $false, "guzman, Carlos Barrera", 2)
